$d = $word.ActiveDocument

# Remove the existing _GoBack bookmark (it will be relocated).
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

# Insert a new paragraph after the first paragraph ("...Project Notes140517")
# containing the text "Dummy line ".
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(2)
$p2.Range.Text = "Dummy line "

# Insert another new, empty paragraph right after it.
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs(3)

# Place a (temporary) placeholder character in the new empty paragraph so we
# can anchor a non-collapsed range there, add the _GoBack bookmark around it,
# then remove the placeholder, leaving a collapsed _GoBack bookmark exactly
# at that empty paragraph.
$p3.Range.Text = "X"
$rngPlaceholder = $d.Range($p3.Range.Start, $p3.Range.Start + 1)
$d.Bookmarks.Add("_GoBack", $rngPlaceholder)

$rngPlaceholder2 = $d.Range($p3.Range.Start, $p3.Range.Start + 1)
$rngPlaceholder2.Text = ""
